$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old D column (style-only cells) that is no longer used
$ws.Range("D1:D8").Clear()

# Write the new dummy data: 17 rows x 3 columns (A, B, C)
for ($r = 1; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = 2400 + $r
    $ws.Cells.Item($r, 2).Value = 121 + ($r - 1) * 2
    $ws.Cells.Item($r, 3).Value = 122 + ($r - 1) * 2
}

# Update the selection to match the target state
$ws.Range("C5").Select()
